$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C15").Value = 'Устраните дипломата, работающего на {TEAM_TAR.FactionDef.Name}'
$ws.Range("C22").Value = '{TEAM_TAR.FactionDef.ShortName} уже несколько месяцев доставляет нам проблемы здесь, в системе {TGT_SYSTEM.name}, и наши переговоры с ними застопорились из-за упрямого дипломата, представителя {TEAM_TAR.FactionDef.Name}. И, предположим, если бы он был, так сказать, «вынесен за скобки», может быть, посредством «случайного» рейда «неизвестных головорезов», мы могли бы добиться некоторого прогресса в дипломатии. Ну как, по рукам?'
$ws.Range("C25").Value = 'Военный преступник, за которым вы охотитесь, будет не один. {TEAM_TAR.FactionDef.Name} выслали отряд наемников, с которыми он подписал контракт. Они начнут атаковать вас без предупреждения.'
$ws.Range("C35").Value = 'Уничтожьте силы, поддерживающие {TEAM_TAR.FactionDef.Name}'
$ws.Range("C45").Value = 'Разыскиваемый военный преступник обнаружен на планете {TGT_SYSTEM.name}. К сожалению, он, похоже, внедрился к наёмникам, работающим на {TEAM_TAR.FactionDef.Name}. Нам нужно, чтобы вы выследили его и уничтожили его мех. Если он катапультируется, мы возьмем его под стражу, чтобы он предстал перед судом {TEAM_EMP.FactionDef.Name}. Если не выживет... ну... оплату вы так или иначе получите.'
$ws.Range("C59").Value = 'Уничтожьте команду засады {TEAM_TAR.FactionDef.Name}'
$ws.Range("C68").Value = 'Мы предлагаем крупную награду за главу местных бандитов. Целью является мятежный мехвоин из {TEAM_TAR.FactionDef.Name}, который атаковал персонал и здания, принадлежащие {TEAM_EMP.FactionDef.Name} в этой системе. Нам нужен отряд наемников, чтобы найти его и уничтожить, прежде чем он перенесёт свою агрессию на наши поселения.'
$ws.Range("C73").Value = 'Миссия выполнена. Пропагандист {TEAM_TAR.FactionDef.Name} мертв.'
$ws.Range("C83").Value = 'Уничтожьте пропагандиста {TEAM_TAR.FactionDef.Name}'
$ws.Range("C90").Value = 'В течение нескольких месяцев у нас были проблемы с людьми здесь, на {TGT_SYSTEM.name}, которых раздражали разговоры о мятеже, восстании и независимости. Основываясь на отчётах ведущих аналитиков {TEAM_EMP.FactionDef.Name}, мы полагаем, что в этом районе есть пропагандист {TEAM_TAR.FactionDef.Name}, вызывающий недовольство. Мы создали защитный периметр, чтобы захватить ее. Но мы не можем покрыть все возможные выходы и могли бы использовать дополнительный отряд мехов. Если она попытается проскользнуть через назначенную вам территорию, превратите ее отряд в шлак.'
$ws.Range("C93").Value = 'Мы наблюдаем много активности от {TEAM_TAR.FactionDef.Name} в непосредственной близости от цели.'
$ws.Range("C113").Value = 'У нас есть возможность устранить ключевую фигуру из {TEAM_TAR.FactionDef.Name} в этой системе, что ввергнет их на время в хаос и обеспечит безопасность операций {TEAM_EMP.FactionDef.Name}. Цель будет относительно беззащитной и уязвимой, но не долго и мы хотели бы, чтобы вы нанесли быстрый удар в это время.'
$ws.Range("C116").Value = 'Мы наблюдаем много активности от {TEAM_TAR.FactionDef.Name} поблизости. Вероятно, это наша цель.'
$ws.Range("C136").Value = 'Наши разведчики полагают, что лидер {TEAM_TAR.FactionDef.Name} присутствует на поверхности {TGT_SYSTEM.name} и уязвим. Его устранение приведёт к беспорядку в их организации и расчистит поле для дальнейших операций {TEAM_EMP.FactionDef.Name}. Цель будет относительно беззащитной и мы хотели бы, чтобы вы нанесли удар, пока есть возможность.'
$ws.Range("C163").Value = 'Наши шпионы сообщают, что ветеран {TEAM_TAR.FactionDef.Name} времён Войн за Наследие ушёл на пенсию и поселился в этой системе. Люди из {TEAM_EMP.FactionDef.ShortName} очень хотели бы, чтобы он и дальше оставался в отставке. И если возможно - навсегда. Мы наблюдаем за десантным кораблем на орбите вокруг планеты и полагаем, что они готовятся подобрать кого-то. Нам нужно, чтобы вы высадились и перехватили этого мехвоина, прежде чем он встретится со своими. Скорее всего, он пилотирует довольно тяжелый мех, поэтому будьте осторожны.'

# Touch column D so its formatting record is materialized separately from
# the rest of the (otherwise identical) trailing column range.
$ws.Columns.Item(4).Hidden = $false | Out-Null

$ws.Range("C10").Select() | Out-Null
